{"js": "// correct some spelling mistakes\nconst body = context.document.body;\n\n// Helper: find the first occurrence of `find` text in the body and replace\n// it with `replace`, case-sensitively and matching whole phrase (not just\n// whole word), then sync.\nasync function replaceFirst(find, replace) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replace, \"Replace\");\n    await context.sync();\n  }\n}\n\n// 1. Fix \"automatisierung\" -> \"Automatisierung\" (capitalize)\nawait replaceFirst(\n  \"durch automatisierung erleichtert.\",\n  \"durch Automatisierung erleichtert.\"\n);\n\n// 2. Insert a comma after \"Es muss m\u00f6glich sein\" in the five bullet points.\nawait replaceFirst(\n  \"Es muss m\u00f6glich sein einen neuen Mitarbeiter anzulegen\",\n  \"Es muss m\u00f6glich sein, einen neuen Mitarbeiter anzulegen\"\n);\nawait replaceFirst(\n  \"Es muss m\u00f6glich sein neue Vorlesungen anzulegen\",\n  \"Es muss m\u00f6glich sein, neue Vorlesungen anzulegen\"\n);\nawait replaceFirst(\n  \"Es muss m\u00f6glich sein R\u00e4ume f\u00fcr Vorlesungen einzuplanen.\",\n  \"Es muss m\u00f6glich sein, R\u00e4ume f\u00fcr Vorlesungen einzuplanen.\"\n);\nawait replaceFirst(\n  \"Es muss m\u00f6glich sein Vorlesungen f\u00fcr Studieng\u00e4nge einzuplanen.\",\n  \"Es muss m\u00f6glich sein, Vorlesungen f\u00fcr Studieng\u00e4nge einzuplanen.\"\n);\nawait replaceFirst(\n  \"Es muss m\u00f6glich sein verschiedene Formen der Ausgabe zu generieren.\",\n  \"Es muss m\u00f6glich sein, verschiedene Formen der Ausgabe zu generieren.\"\n);\n\n// 3. \"einen View\" -> \"ein View\"\nawait replaceFirst(\n  \"Zus\u00e4tzlich ben\u00f6tigt er einen View mit Detailinformationen\",\n  \"Zus\u00e4tzlich ben\u00f6tigt er ein View mit Detailinformationen\"\n);\n\n// 4. Move the _GoBack bookmark from the first paragraph to the end of the\n//    \"... ein View mit Detailinformationen einer jeden Veranstaltung. \"\n//    paragraph (content-only span, i.e. before the paragraph mark).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Zus\u00e4tzlich ben\u00f6tigt er ein View mit Detailinformationen\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const endRange = target.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# correct some spelling mistakes\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1. Fix \"automatisierung\" -> \"Automatisierung\" (capitalize)\nReplace-Text \"durch automatisierung erleichtert.\" \"durch Automatisierung erleichtert.\"\n\n# 2. Insert comma after \"Es muss m\u00f6glich sein\" in the 5 bullet points\nReplace-Text \"Es muss m\u00f6glich sein einen neuen Mitarbeiter anzulegen\" \"Es muss m\u00f6glich sein, einen neuen Mitarbeiter anzulegen\"\nReplace-Text \"Es muss m\u00f6glich sein neue Vorlesungen anzulegen\" \"Es muss m\u00f6glich sein, neue Vorlesungen anzulegen\"\nReplace-Text \"Es muss m\u00f6glich sein R\u00e4ume f\u00fcr Vorlesungen einzuplanen.\" \"Es muss m\u00f6glich sein, R\u00e4ume f\u00fcr Vorlesungen einzuplanen.\"\nReplace-Text \"Es muss m\u00f6glich sein Vorlesungen f\u00fcr Studieng\u00e4nge einzuplanen.\" \"Es muss m\u00f6glich sein, Vorlesungen f\u00fcr Studieng\u00e4nge einzuplanen.\"\nReplace-Text \"Es muss m\u00f6glich sein verschiedene Formen der Ausgabe zu generieren.\" \"Es muss m\u00f6glich sein, verschiedene Formen der Ausgabe zu generieren.\"\n\n# 3. \"einen View\" -> \"ein View\"\nReplace-Text \"Zus\u00e4tzlich ben\u00f6tigt er einen View mit Detailinformationen\" \"Zus\u00e4tzlich ben\u00f6tigt er ein View mit Detailinformationen\"\n\n# 4. Move the _GoBack bookmark from the first paragraph to the end of the\n#    \"... ein View mit Detailinformationen einer jeden Veranstaltung. \" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"ein View mit Detailinformationen einer jeden Veranstaltung.\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$found = $find2.Execute()\nif ($found) {\n    $target = $d.Content\n    $target.Start = $find2.Parent.End\n    $target.End = $find2.Parent.End\n    $d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n}\n"}
